$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: row, NACE code (numeric), Name, AG
$data = @(
    @(2, 4110, "F-4110 Development of building projects", "Y"),
    @(3, 4211, "F-4211 Construction of roads and motorways", "G"),
    @(4, 4120, "F-4120 Construction of residential and non-residential buildings", "Z"),
    @(5, 4221, "F-4221 Construction of utility projects for fluids", "F"),
    @(6, 4399, "F-4399 Other specialised construction activities n.e.c.", "F"),
    @(7, 4391, "F-4391 Roofing activities", "F"),
    @(8, 4222, "F-4222 Construction of utility projects for electricity and telecommunications", "F"),
    @(9, 4291, "F-4291 Construction of water projects", "F"),
    @(10, 4299, "F-4299 Construction of other civil engineering projects n.e.c.", "F"),
    @(11, 4212, "F-4212 Construction of railways and underground railwaysÊ", "F"),
    @(12, 4333, "F-4333 Floor and wall covering", "F"),
    @(13, 4329, "F-4329 Other construction installation", "F"),
    @(14, 4332, "F-4332 Joinery installation", "F"),
    @(15, 4213, "F-4213 Construction of bridges and tunnels", "F"),
    @(16, 4339, "F-4339 Other building completion and finishing", "F"),
    @(17, 4312, "F-4312 Site preparation", "F"),
    @(18, 4311, "F-4311 Demolition", "F"),
    @(19, 3821, "E-3821 Treatment and disposal of non-hazardous waste", "E"),
    @(20, 3811, "E-3811 Collection of non-hazardous waste", "E"),
    @(21, 3832, "E-3832 Recovery of sorted materials", "E"),
    @(22, 3900, "E-3900 Remediation activities and other waste management services", "E"),
    @(23, 3822, "E-3822 Treatment and disposal of hazardous waste", "E"),
    @(24, 3812, "E-3812 Collection of hazardous waste", "E")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("D13").Select()
